$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 376.25
$ws.Range("I28").Value = 368.33334
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 368.33334
$ws.Range("L28").Value = 400
$ws.Range("M28").Value = 116.66666
$ws.Range("N28").Value = -1370
$ws.Range("H40").Value = 1799.4615
$ws.Range("I40").Value = 1099.2222
$ws.Range("J40").Value = 3375
$ws.Range("K40").Value = 1099.2222
$ws.Range("L40").Value = 3375
$ws.Range("M40").Value = -924.2221999999999
$ws.Range("N40").Value = -3725
$ws.Range("H64").Value = 4448.636
$ws.Range("I64").Value = 3666.3333
$ws.Range("J64").Value = 4526.8667
$ws.Range("K64").Value = 3666.3333
$ws.Range("L64").Value = 4526.8667
$ws.Range("M64").Value = -3418.3333
$ws.Range("N64").Value = -5022.8667
$ws.Range("H67").Value = 4448.636
$ws.Range("I67").Value = 3666.3333
$ws.Range("J67").Value = 4526.8667
$ws.Range("K67").Value = 3666.3333
$ws.Range("L67").Value = 4526.8667
$ws.Range("M67").Value = -2808.3333
$ws.Range("N67").Value = -6242.8667
$ws.Range("H76").Value = 4202.615
$ws.Range("I76").Value = 4104.857
$ws.Range("J76").Value = 4316.6665
$ws.Range("K76").Value = 4104.857
$ws.Range("L76").Value = 4316.6665
$ws.Range("M76").Value = -3789.857
$ws.Range("N76").Value = -4946.6665
$ws.Range("H79").Value = 4202.615
$ws.Range("I79").Value = 4104.857
$ws.Range("J79").Value = 4316.6665
$ws.Range("K79").Value = 4104.857
$ws.Range("L79").Value = 4316.6665
$ws.Range("M79").Value = -3012.857
$ws.Range("N79").Value = -6500.6665
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1025
$ws.Range("K107").Value = 1025
$ws.Range("M107").Value = 895
$ws.Range("H129").Value = 970.4643
$ws.Range("J129").Value = 1440.3572
$ws.Range("L129").Value = 4321.071599999999
$ws.Range("N129").Value = -14321.0716
$ws.Range("H140").Value = 108317.14
$ws.Range("J140").Value = 110778.336
$ws.Range("L140").Value = 110778.336
$ws.Range("N140").Value = -121138.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1243.6666
$ws.Range("I45").Value = 1278.6
$ws.Range("K45").Value = 1278.6
$ws.Range("M45").Value = -901.5999999999999
$ws.Range("H102").Value = 3300
$ws.Range("I102").Value = 1525
$ws.Range("K102").Value = 1525
$ws.Range("M102").Value = 97
$ws.Range("H132").Value = 716749.6
$ws.Range("I132").Value = 911017.8
$ws.Range("K132").Value = 2733053.4
$ws.Range("M132").Value = -2730523.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3418.7856
$ws.Range("I105").Value = 3488.5833
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 3488.5833
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -1741.5833
$ws.Range("N105").Value = -6494
$ws.Range("H134").Value = 515586.7
$ws.Range("J134").Value = 4181.4
$ws.Range("L134").Value = 12544.2
$ws.Range("N134").Value = -17614.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H62").Value = 85842.336
$ws.Range("I62").Value = 168784.67
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 168784.67
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -168160.67
$ws.Range("N62").Value = -4148
$ws.Range("H65").Value = 85842.336
$ws.Range("I65").Value = 168784.67
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 843923.3500000001
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -840803.3500000001
$ws.Range("N65").Value = -20740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1073586.4
$ws.Range("I12").Value = 38.57143
$ws.Range("J12").Value = 1756753.1
$ws.Range("K12").Value = 115.71429
$ws.Range("L12").Value = 5270259.300000001
$ws.Range("M12").Value = 57.28570999999999
$ws.Range("N12").Value = -5270605.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5734.613
$ws.Range("I70").Value = 5001.2666
$ws.Range("J70").Value = 6422.125
$ws.Range("K70").Value = 5001.2666
$ws.Range("L70").Value = 6422.125
$ws.Range("M70").Value = -4731.2666
$ws.Range("N70").Value = -6962.125
$ws.Range("H73").Value = 5734.613
$ws.Range("I73").Value = 5001.2666
$ws.Range("J73").Value = 6422.125
$ws.Range("K73").Value = 5001.2666
$ws.Range("L73").Value = 6422.125
$ws.Range("M73").Value = -4065.2666
$ws.Range("N73").Value = -8294.125
$ws.Range("H80").Value = 3467.2222
$ws.Range("I80").Value = 3367.5
$ws.Range("J80").Value = 3666.6667
$ws.Range("K80").Value = 3367.5
$ws.Range("L80").Value = 3666.6667
$ws.Range("M80").Value = -2369.5
$ws.Range("N80").Value = -5662.6667
$ws.Range("H83").Value = 3467.2222
$ws.Range("I83").Value = 3367.5
$ws.Range("J83").Value = 3666.6667
$ws.Range("K83").Value = 16837.5
$ws.Range("L83").Value = 18333.3335
$ws.Range("M83").Value = -11845.5
$ws.Range("N83").Value = -28317.3335
$ws.Range("H93").Value = 20250
$ws.Range("J93").Value = 20250
$ws.Range("L93").Value = 20250
$ws.Range("N93").Value = -23994
$ws.Range("H102").Value = 2029.5454
$ws.Range("I102").Value = 2014.7646
$ws.Range("J102").Value = 2079.8
$ws.Range("K102").Value = 2014.7646
$ws.Range("L102").Value = 2079.8
$ws.Range("M102").Value = -392.7646
$ws.Range("N102").Value = -5323.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 14712140
$ws.Range("J122").Value = 5980.8
$ws.Range("L122").Value = 17942.4
$ws.Range("N122").Value = -22842.4
$ws.Range("I132").Value = 3309.875
$ws.Range("J132").Value = 5599.7144
$ws.Range("K132").Value = 9929.625
$ws.Range("L132").Value = 16799.1432
$ws.Range("M132").Value = -7399.625
$ws.Range("N132").Value = -21859.1432
$ws.Range("H135").Value = 135588.89
$ws.Range("J135").Value = 135588.89
$ws.Range("L135").Value = 135588.89
$ws.Range("N135").Value = -145728.89

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 32500
$ws.Range("H113").Value = 865.5714
$ws.Range("I113").Value = 951.2
$ws.Range("J113").Value = 651.5
$ws.Range("K113").Value = 2853.6
$ws.Range("L113").Value = 1954.5
$ws.Range("M113").Value = -683.6000000000004
$ws.Range("N113").Value = -6294.5
$ws.Range("H123").Value = 23966.924
$ws.Range("J123").Value = 23966.924
$ws.Range("L123").Value = 23966.924
$ws.Range("N123").Value = -33766.924
$ws.Range("H132").Value = 1370.1143
$ws.Range("I132").Value = 1244.1072
$ws.Range("J132").Value = 1874.1428
$ws.Range("K132").Value = 3732.3216
$ws.Range("L132").Value = 5622.428400000001
$ws.Range("M132").Value = -1202.3216
$ws.Range("N132").Value = -10682.4284
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120
